$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "row 2" and "row 3" data for the fields that differ between
# the two price-record rows (date, volume, prices, origin, price/kg).

# Row 2 original values (becomes row 3 target)
$d2 = $ws.Range("D2").Value2
$m2 = $ws.Range("M2").Value2
$n2 = $ws.Range("N2").Value2
$o2 = $ws.Range("O2").Value2
$p2 = $ws.Range("P2").Value2
$r2 = $ws.Range("R2").Value2
$s2 = $ws.Range("S2").Value2

# Row 3 original values (becomes row 2 target)
$d3 = $ws.Range("D3").Value2
$m3 = $ws.Range("M3").Value2
$n3 = $ws.Range("N3").Value2
$o3 = $ws.Range("O3").Value2
$p3 = $ws.Range("P3").Value2
$r3 = $ws.Range("R3").Value2
$s3 = $ws.Range("S3").Value2

# Write swapped values into row 2
$ws.Range("D2").Value = $d3
$ws.Range("M2").Value = $m3
$ws.Range("N2").Value = $n3
$ws.Range("O2").Value = $o3
$ws.Range("P2").Value = $p3
$ws.Range("R2").Value = $r3
$ws.Range("S2").Value = $s3

# Write swapped values into row 3
$ws.Range("D3").Value = $d2
$ws.Range("M3").Value = $m2
$ws.Range("N3").Value = $n2
$ws.Range("O3").Value = $o2
$ws.Range("P3").Value = $p2
$ws.Range("R3").Value = $r2
$ws.Range("S3").Value = $s2
